# "removed path prefix from ms method"
#
# 1. Strip the "D:\Projects\Default\Methods\" directory prefix from the two
#    6560_Method shared strings (POS / NEG method file names).
# 2. Cosmetic view-state changes that came along with the edit:
#    - samples sheet becomes the active/selected tab, zoomed to 231%,
#      with the selection parked on D10.
#    - rf_params sheet is no longer the selected tab; its lingering
#      selection moves to I18.
#    - samples header row (A1:J1) becomes bold.
#    - a handful of samples columns get explicit widths.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # samples
$ws2 = $wb.Worksheets.Item(2)   # rf_params

# --- 1. Remove the hard-coded path prefix from the method file names ---
$ws1.Range("A1:J123").Replace("D:\Projects\Default\Methods\", "")

# --- 2. Column widths on the samples sheet ---
$ws1.Columns.Item(2).ColumnWidth = 12.5               # B -> 13.33203125
$ws1.Columns.Item(4).ColumnWidth = 14.333333333333334 # D -> 15.1640625
$ws1.Columns.Item(5).ColumnWidth = 16                 # E -> 16.83203125
$ws1.Columns.Item(6).ColumnWidth = 12.333333333333334 # F -> 13.1640625
$ws1.Columns.Item(7).ColumnWidth = 12.833333333333334 # G -> 13.6640625
$ws1.Columns.Item(9).ColumnWidth = 11.5               # I -> 12.33203125

# --- 3. Bold the samples header row ---
$ws1.Range("A1:J1").Font.Bold = $true

# --- 4. View / selection state ---
# Leave rf_params first so the samples tab ends up as the final active one.
$ws2.Activate()
$ws2.Range("I18").Select()

$ws1.Activate()
$excel.ActiveWindow.Zoom = 231
$ws1.Range("D10").Select()
